# Update contact list: replace Pierre Vanobbergen's row with a second
# Melvin Leble entry (capitalized surname + new personal email), and
# update the active cell selection, per the "Add Inventory Manager for
# serialnumber" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 stays Melvin Leble / melvin.leble@supinfo.com (re-written, same values)
$ws.Range("A2").Value = "Melvin"
$ws.Range("B2").Value = "Leble"
$ws.Range("C2").Value = "melvin.leble@supinfo.com"

# Row 3 used to be Pierre Vanobbergen; now it's Melvin LEBLE with a personal email
$ws.Range("A3").Value = "Melvin"
$ws.Range("B3").Value = "LEBLE"
$ws.Range("C3").Value = "leble17@gmail.com"

# Move the active selection to A4, matching the saved view state
[void]$ws.Range("A4").Select()
